$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 56 (pushing the existing
# rows 56-105 down to 57-106). Insert a fresh row and populate it.
$ws.Rows.Item(56).Insert()

$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "Vega Modelo de Temuco"
$ws.Range("C56").Value = "La Araucanía"
$ws.Range("D56").Value = 44566
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = 100112031
$ws.Range("G56").Value = "Poroto verde"
$ws.Range("H56").Value = "Brío"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 80
$ws.Range("K56").Value = 2000
$ws.Range("L56").Value = 2000
$ws.Range("M56").Value = 2000
$ws.Range("N56").Value = "$/kilo"
$ws.Range("O56").Value = "Región de La Araucanía"
$ws.Range("P56").Value = 2000
$ws.Range("Q56").Value = 1
$ws.Range("R56").Value = "Hortaliza"
